# ---------------------------------------------------------------------------
# Commit: "Wed, Apr 15, 2020  5:06:23 PM"
#
# The underlying OOXML edit swaps the contents of ppt/theme/theme1.xml and
# ppt/theme/theme2.xml:
#   - theme1.xml (was "Office Theme" / "Office" colors, used only by the
#     notes master) becomes the "Integral" / "Red Violet" theme.
#   - theme2.xml (was "Integral" / "Red Violet", the theme actually used by
#     the slide master / the visible design) becomes "Office Theme" /
#     "Office" colors.
#
# In other words: the colour palette that PowerPoint actually shows on the
# slides switches from the pink/violet "Integral" palette to the default
# "Office" palette. The font scheme and format (fill/line/effect) scheme are
# byte-for-byte identical between the two theme parts, so the only visible
# change is the 12 theme colours (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# We reach that palette through the Design/Master object model, which is the
# supported COM surface for editing a theme's colour scheme.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# The presentation's one-and-only design/theme (as used by the slide master)
# backs ppt/theme/theme2.xml.
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Target palette: the "Office Theme" colours (standard PowerPoint defaults).
# Index order follows the ThemeColorScheme.Colors COM indexing:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeTheme = @{
    1  = 0x000000   # dk1
    2  = 0xFFFFFF   # lt1
    3  = 0x44546A   # dk2
    4  = 0xE7E6E6   # lt2
    5  = 0x5B9BD5   # accent1
    6  = 0xED7D31   # accent2
    7  = 0xA5A5A5   # accent3
    8  = 0xFFC000   # accent4
    9  = 0x4472C4   # accent5
    10 = 0x70AD47   # accent6
    11 = 0x0563C1   # hlink
    12 = 0x954F72   # folHlink
}

function ConvertTo-OleRgb($rrggbb) {
    # PowerPoint's RGB() values are packed little-endian (0x00BBGGRR) while
    # the hex values above are plain 0xRRGGBB, so swap R and B.
    $r = ($rrggbb -shr 16) -band 0xFF
    $g = ($rrggbb -shr 8) -band 0xFF
    $b = $rrggbb -band 0xFF
    return ($b -shl 16) -bor ($g -shl 8) -bor $r
}

for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = ConvertTo-OleRgb $officeTheme[$i]
}

Write-Host "Theme colour scheme updated to Office Theme palette."
